$wb = $excel.ActiveWorkbook

# --- Pending_Orders sheet ---
$po = $wb.Worksheets.Item("Pending_Orders")

# Update the limit price for the CB order (row 4, column C)
$po.Cells.Item(4, 3).Value = 280

# Remove the last 5 rows (rows 7-11), which were duplicates of rows 2-6
$po.Range("A7:I11").Delete()

# --- Benchmark sheet ---
$bm = $wb.Worksheets.Item("Benchmark")

# Append two more rows (6 and 7) that duplicate the existing row 5 data.
# Format column A as text first so the date-like string isn't
# auto-converted into a date serial number, then restore the default
# "Normal" style so the cell matches the plain (unstyled) string cells
# already used by rows 2-5 in column A.
$dateRange = $bm.Range("A6:A7")
$dateRange.NumberFormat = "@"
$dateRange.Value = "2026-01-26"
$dateRange.NumberFormat = "General"
$dateRange.Style = "Normal"

foreach ($r in 6, 7) {
    $bm.Cells.Item($r, 2).Value = 637.0900268554688
    $bm.Cells.Item($r, 3).Value = 44876.20025634766
    $bm.Cells.Item($r, 4).Value = 0.001572146796001695
    $bm.Cells.Item($r, 5).Value = 0.1219050064086914
    $bm.Cells.Item($r, 6).Value = 0.1203328596126897
}

# Carry over the same number format as row 5 for the percentage columns (D, E, F)
$bm.Range("D5:F5").Copy() | Out-Null
$bm.Range("D6:F7").PasteSpecial(-4122) | Out-Null
